# Switch in designlabels switched to numeric
# Column B ("S" = Switch) held text labels "power_on"/"power_off".
# Replace them with their numeric equivalents: power_on -> 1, power_off -> 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 2)  # column B
    $val = $cell.Value2
    if ($val -eq "power_on") {
        $cell.Value = 1
    } elseif ($val -eq "power_off") {
        $cell.Value = 0
    }
}

# Update the active selection to match the recorded view state (B17).
$ws.Range("B17").Select()
